$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.477.26'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.572.77'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '291.92'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3722'
$ws.Range('E7').Value = '  -1.23%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '49.88'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3401'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.25'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.046'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.967'
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '1.571.51'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001125'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.79'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06761'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.308'
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.37'
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('D24').Value = '22.469.08'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.373'
$ws.Range('E25').Value = '  -0.91%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.628'
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.03'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '149.48'
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.047'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.38'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('D31').Value = '1.746.52'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.085'
$ws.Range('E32').Value = '  +9.71%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.207'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.014'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.810'
$ws.Range('E35').Value = '  -3.58%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08357'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02481'
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2301'
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('E39').Value = '  -3.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06546'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.448'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6239'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '14.05'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.816'
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5847'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '130.78'
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.073'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.212'
$ws.Range('E50').Value = '  -5.44%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07334'
$ws.Range('E51').Value = '  +0.05%  '
